$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the data rows (2-5, columns A-F) with the new template value used for
# generating the experiment plan, matching row 1's height (19.5) and clearing
# the #,##0 number format down to General (border/font/alignment stay as-is).
$ws.Range("A1:F5").NumberFormat = "general"

$ws.Rows("2:5").RowHeight = 19.5

$ws.Range("A2:F5").Value = "1 4 45"
$ws.Range("A2").Value = "0 0 0"
$ws.Range("E5").Value = "0 0 0"
